# laba2gribach.docx edits:
# 1) "Грибач Н.Э." paragraph: add left indent (5760 twips = 288 pt) and
#    change justification from right to center.
# 2) "Проверил:" paragraph: the last (9th) tab run becomes a bookmarked
#    run with the text "Усенко Ф.В.".
# 3) Embedded OLEObject's ObjectID attribute changes.

$d = $word.ActiveDocument

# --- Edit 1: paragraph formatting for "Грибач Н.Э." ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Грибач Н\.Э\.") {
        $p.Format.LeftIndent = 288
        $p.Format.Alignment = 1
        break
    }
}

# --- Edit 2: replace last tab in the "Проверил:" paragraph ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^Проверил:") {
        $pEnd = $p.Range.End
        $tabEnd = $pEnd - 1
        $tabStart = $tabEnd - 1
        $tabRange = $d.Range($tabStart, $tabEnd)
        $tabRange.Text = "Усенко Ф.В."
        $tabRange.Font.Name = "Times New Roman"
        $tabRange.Font.Size = 14
        $tabRange.LanguageID = 1049
        $d.Bookmarks.Add("_Hlk183354226", $tabRange)
        break
    }
}

# --- Edit 3: OLEObject's ObjectID ---
$d.Content.Find.Execute("_1793968453", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "_1794604514", 2) | Out-Null
